# Update "To Do.xlsx" / Sheet1: arduino file updated to use level switches
# so ball valves shut off when a level switch is activated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Move the "Date" value back one day (30948 -> 30947)
$ws.Range("A2").Value = (Get-Date -Year 1984 -Month 9 -Day 22 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)

# --- Add the solution note for the level-switch integration row
$ws.Range("C3").Value = "Integrated, 9/23/19"

# --- Header row: fix "Solution\" -> "Solution"
$ws.Range("C1").Value = "Solution"

# --- Grow the merged "Date" cell down so it spans rows 2-6 (previously 2-4),
#     matching the level-switch/ball-valve/arduino rows underneath it.
$ws.Range("A2:A4").MergeCells = $false
$ws.Range("A2:A6").Merge()

# --- Widen column C to fit the new "Integrated, 9/23/19" text
$ws.Columns.Item(3).ColumnWidth = 17.5

# --- Restore the selection as recorded after the edit
$ws.Range("B18").Select() | Out-Null
